$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing "sample=0.01" rows that previously lived at rows 18-19
# (the new curve only goes down to 0.01 once, at rows 16-17) and shift every
# remaining row's "sample" value up into the slot above it.
$ws.Rows("18:19").Delete()

# Rewrite column A (sample), and populate the new C (best_dilution) / D
# (plate_count) columns for every data row. Column B (rep: a/b) already
# alternates correctly and needs no changes.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 3).Value = 1000000
$ws.Cells.Item(2, 4).Value = 131

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 3).Value = 1000000
$ws.Cells.Item(3, 4).Value = 94

$ws.Cells.Item(4, 1).Value = 0.75
$ws.Cells.Item(4, 3).Value = 1000000
$ws.Cells.Item(4, 4).Value = 45

$ws.Cells.Item(5, 1).Value = 0.75
$ws.Cells.Item(5, 3).Value = 1000000
$ws.Cells.Item(5, 4).Value = 46

$ws.Cells.Item(6, 1).Value = 0.5
$ws.Cells.Item(6, 3).Value = 100000
$ws.Cells.Item(6, 4).Value = 175

$ws.Cells.Item(7, 1).Value = 0.5
$ws.Cells.Item(7, 3).Value = 100000
$ws.Cells.Item(7, 4).Value = 228

$ws.Cells.Item(8, 1).Value = 0.25
$ws.Cells.Item(8, 3).Value = 100000
$ws.Cells.Item(8, 4).Value = 173

$ws.Cells.Item(9, 1).Value = 0.25
$ws.Cells.Item(9, 3).Value = 100000
$ws.Cells.Item(9, 4).Value = 160

$ws.Cells.Item(10, 1).Value = 0.1
$ws.Cells.Item(10, 3).Value = 100000
$ws.Cells.Item(10, 4).Value = 67

$ws.Cells.Item(11, 1).Value = 0.1
$ws.Cells.Item(11, 3).Value = 100000
$ws.Cells.Item(11, 4).Value = 95

$ws.Cells.Item(12, 1).Value = 0.05
$ws.Cells.Item(12, 3).Value = 100000
$ws.Cells.Item(12, 4).Value = 36

$ws.Cells.Item(13, 1).Value = 0.05
$ws.Cells.Item(13, 3).Value = 100000
$ws.Cells.Item(13, 4).Value = 89

$ws.Cells.Item(14, 1).Value = 0.025
$ws.Cells.Item(14, 3).Value = 10000
$ws.Cells.Item(14, 4).Value = 146

$ws.Cells.Item(15, 1).Value = 0.025
$ws.Cells.Item(15, 3).Value = 10000
$ws.Cells.Item(15, 4).Value = 125

$ws.Cells.Item(16, 1).Value = 0.01
$ws.Cells.Item(16, 3).Value = 10000
$ws.Cells.Item(16, 4).Value = 2

$ws.Cells.Item(17, 1).Value = 0.01
$ws.Cells.Item(17, 3).Value = 10000
$ws.Cells.Item(17, 4).Value = 52

# The new best_dilution column is displayed in scientific notation.
$ws.Range("C2:C17").NumberFormat = "0.00E+00"

# Re-apply the descending sort on column A over the new (smaller) range so
# the worksheet's stored sortState matches the shrunk data extent.
$sortRange = $ws.Range("A2:D17")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A17"), 0, 2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Apply()

# Move the active selection to match the author's final cursor position.
$ws.Range("F15").Select()
